$wb = $excel.ActiveWorkbook

# --- Logs sheet: append new row 24 ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A24").Value = "CE-certificaten verzoek"
$ws.Range("B24").Value = "inkoop@testbedrijf123.nl"
$ws.Range("C24").Value = "Kun je mij de CE-certificaten van de EcoPro-700 sturen?"
$ws.Range("D24").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E24").Value = "Bedankt, we hebben dit doorgestuurd naar kwaliteit@testbedrijf123.nl."
$ws.Range("F24").Value = "2025-08-14 21:20:14"
$ws.Range("G24").Value = "Nee"
$ws.Range("H24").Value = "Ja"
$ws.Range("I24").Value = "Nee"
$ws.Range("J24").Value = "Nee"

# --- Extend the conditional-formatting ranges to cover the new row ---
$cfColumns = @("D", "G", "H", "I", "J")
foreach ($col in $cfColumns) {
    $oldRange = $col + "2:" + $col + "23"
    $newRange = $col + "2:" + $col + "24"
    $fcs = $ws.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($newRange))
    }
}

# --- Dashboard sheet: bump the "Intern verzoek / Actie voor medewerker" count ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 18
